$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal TEXT value into a cell without Excel's automatic
# "this looks like a number" coercion (e.g. "515150" -> 515150). We do this
# by putting a text-returning formula in a scratch cell, copying it, and
# pasting-special-VALUES into the destination; the pasted result is a plain
# text value (no formula, no residual formatting) in the destination cell.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($scratch, $targetRange, [string]$text)

    $escaped = $text.Replace("""", """""")
    $scratch.Formula = "=""" + $escaped + """"
    $scratch.Copy()
    $targetRange.PasteSpecial(-4163)
}

# scratch cell, far away from any real data, cleared at the very end
$scratchSheet = $wb.Worksheets.Item(1)
$scratch = $scratchSheet.Range("Z100")

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right before "总计".
#    NOTE: sheet handles in this host are position-based, so after the insert
#    shifts "总计" one slot to the right we MUST re-fetch it by name - the
#    old $totalSheet variable would otherwise silently keep pointing at
#    whatever sheet now sits in that old slot (i.e. the new "2022-Q1" sheet).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Add($totalSheet)
$q1Sheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item("总计")

# Clone layout + styling from "2021-Q4" (identical header wording / formats),
# then overwrite with 2022-Q1's own figures. A1 is intentionally skipped (kept
# untouched/empty) since the source sheets never populate that corner cell.
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Range("B1:H4").Copy()
$q1Sheet.Range("B1:H4").PasteSpecial(-4122)
$q1Sheet.Range("B1:H4").PasteSpecial(-4163)
$q4Sheet.Range("A2:A4").Copy()
$q1Sheet.Range("A2:A4").PasteSpecial(-4122)
$q1Sheet.Range("A2:A4").PasteSpecial(-4163)

# Row 2 - 515150 富国中证国企一带一路ETF
$q1Sheet.Range("A2").Value = 0
Set-TextValue $scratch $q1Sheet.Range("B2") "515150"
Set-TextValue $scratch $q1Sheet.Range("C2") "富国中证国企一带一路ETF"
Set-TextValue $scratch $q1Sheet.Range("D2") "7.41"
Set-TextValue $scratch $q1Sheet.Range("E2") "99.21"
Set-TextValue $scratch $q1Sheet.Range("F2") "2.09"
Set-TextValue $scratch $q1Sheet.Range("G2") "0.1549"
$q1Sheet.Range("H2").Value = 10

# Row 3 - 515110 易方达中证国企一带一路ETF
$q1Sheet.Range("A3").Value = 1
Set-TextValue $scratch $q1Sheet.Range("B3") "515110"
Set-TextValue $scratch $q1Sheet.Range("C3") "易方达中证国企一带一路ETF"
Set-TextValue $scratch $q1Sheet.Range("D3") "4.83"
Set-TextValue $scratch $q1Sheet.Range("E3") "99.52"
Set-TextValue $scratch $q1Sheet.Range("F3") "2.09"
Set-TextValue $scratch $q1Sheet.Range("G3") "0.1009"
$q1Sheet.Range("H3").Value = 10

# Row 4 - 515990 汇添富中证国企一带一路ETF
$q1Sheet.Range("A4").Value = 2
Set-TextValue $scratch $q1Sheet.Range("B4") "515990"
Set-TextValue $scratch $q1Sheet.Range("C4") "汇添富中证国企一带一路ETF"
Set-TextValue $scratch $q1Sheet.Range("D4") "1.08"
Set-TextValue $scratch $q1Sheet.Range("E4") "99.16"
Set-TextValue $scratch $q1Sheet.Range("F4") "2.10"
Set-TextValue $scratch $q1Sheet.Range("G4") "0.0227"
$q1Sheet.Range("H4").Value = 10

# ---------------------------------------------------------------------------
# 2. Update "总计" with a new leading row for 2022-Q1, pushing the older
#    quarters down (same 3-column table, now 4 rows of data instead of 3).
# ---------------------------------------------------------------------------

# Give the brand-new A4 index cell the same styling as the existing A column
# cells (A2/A3) before writing into it.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

# Shift 2021-Q4 and 2021-Q3 rows down one position first ...
$totalSheet.Range("A4").Value = 2
Set-TextValue $scratch $totalSheet.Range("B4") "2021-Q3"
$totalSheet.Range("C4").Value = 6
$totalSheet.Range("D4").Value = 0.11

$totalSheet.Range("A3").Value = 1
Set-TextValue $scratch $totalSheet.Range("B3") "2021-Q4"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 0.78

# ... then write the brand-new 2022-Q1 summary row on top.
$totalSheet.Range("A2").Value = 0
Set-TextValue $scratch $totalSheet.Range("B2") "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.28

# ---------------------------------------------------------------------------
# 3. Clean up the scratch cell used for text coercion.
# ---------------------------------------------------------------------------
$scratch.ClearContents()
